# Fixed LHS sampling to only sample across uncertainties (X) that vary
# (Ls still vary for all Ls) and rebuilt templates with
# PFLO:ALL_NO_STOPPING_DEFORESTATION_PLUR
#
# Concretely: the old "strategy_id-5008" template sheet becomes
# "strategy_id-5007", and a brand-new "strategy_id-5009" sheet (an exact
# duplicate of that template) is added right after it.

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("strategy_id-5008")

# Duplicate the template sheet, inserting the copy immediately after it,
# while it still carries its original name.
$srcSheet.Copy($null, $srcSheet)

# Rename the original sheet to "strategy_id-5007" ...
$srcSheet.Name = "strategy_id-5007"

# ... and the newly-inserted duplicate to "strategy_id-5009".
$newSheet = $wb.Worksheets.Item("strategy_id-5008 (2)")
$newSheet.Name = "strategy_id-5009"
